# Apply edits described by the commit "updated spreadsheet with insights from evidently"
# - add a new row (7) with data about the Honeycomb.io evidently case
# - adjust column widths / row heights
# - freeze the header row
# - add (orphaned) conditional-format differential styles (dxfs) to styles.xml

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New row 7 content (shared strings 28-33 in the target workbook)
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 'https://www.honeycomb.io/blog/hard-stuff-nobody-talks-about-llm'
$ws.Range("B7").Value = 'Honeycomb.io
Тула для аналитики работы сервисов и продуктов'
$ws.Range("C7").Value = 'Фактически делали text2sql :
- пихали на вход гпт схему данных, текстовый запрос юзера, и несколько примеров
- пытались обойти лимит токенов через итерационный вызов ЛЛМ и эмбеддинги
- экспериментировали с промптами, чтобы добиться устойчивости и скорости работы '
$ws.Range("D7").Value = 'не говорят явно
в статье только субъективные оценки
+
процент получившихся валидных запросов (sql-like)'
$ws.Range("E7").Value = 'не говорят явно'
$ws.Range("F7").Value = 'тюнили промпты'

$ws.Rows.Item(7).RowHeight = 80.55

# ---------------------------------------------------------------------------
# 2. Row height tweaks for existing rows
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 70.6
$ws.Rows.Item(6).RowHeight = 91.5

# ---------------------------------------------------------------------------
# 3. Column width tweaks
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 13
$ws.Columns.Item(3).ColumnWidth = 56.166666666666664
$ws.Columns.Item(6).ColumnWidth = 25.833333333333332

# ---------------------------------------------------------------------------
# 4. Freeze header row (row 1) and restore view/selection state
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F8").Select()

# ---------------------------------------------------------------------------
# 5. Orphaned conditional-format dxfs (white + blue fills), added then the
#    rule itself removed again - mirrors what is observed in the target file
#    (dxfs present in styles.xml but no conditionalFormatting left on sheet)
# ---------------------------------------------------------------------------
$fc = $ws.Range("A1").FormatConditions
$cond1 = $fc.Add(1, 3, "5")
$cond1.Interior.Color = 16777215
$cond2 = $fc.Add(1, 3, "5")
$cond2.Interior.Color = 16711680
$fc.Delete()

Write-Host "edits applied"
